$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("D4").Value = "MCT-1A-Circuitos elétricos"
$ws.Range("E4").Value = "-"

# Row 6
$ws.Range("D6").Value = "MCT-1A-Circuitos elétricos"
$ws.Range("E6").Value = "MCT-3A-Máquinas Elétricas"
$ws.Range("F6").Value = "MCT-1A-Circuitos elétricos"

# Row 7
$ws.Range("E7").Value = "MCT-3A-Máquinas Elétricas"
$ws.Range("F7").Value = "MCT-1A-Circuitos elétricos"

# Row 20
$ws.Range("D20").Value = "['ELM-2NA-Automação Industrial', 'ELM-2NA-Automação Industrial', 'ELM-2NA-Automação Industrial', 'ELM-2NA-Automação Industrial']"

# Row 21
$ws.Range("C21").Value = "['ELM-1NA-Sistemas Digitais', 'ELM-1NA-Sistemas Digitais', 'ELM-1NA-Sistemas Digitais', 'ELM-1NA-Sistemas Digitais']"
